# 2019-12-24 공부한 내용 추가
# Mark "실제 결과" (E) / "테스트 결과" (G) columns with "O" for the test rows
# that already have a tester ("Test 담당자" / F column) filled in, on the
# 메인화면, 로그인&로그아웃 and 회원가입 sheets. Also update the active
# sheet/selection state to match the saved workbook (메인화면 becomes the
# active tab).

$wb = $excel.ActiveWorkbook

# --- 메인화면 (sheet 2): rows 11-20 ---
$wsMain = $wb.Worksheets.Item(2)
$wsMain.Range("E11:E20").Value = "O"
$wsMain.Range("G11:G20").Value = "O"

# --- 로그인&로그아웃 (sheet 3): rows 12-15 ---
$wsLogin = $wb.Worksheets.Item(3)
$wsLogin.Range("E12:E15").Value = "O"
$wsLogin.Range("G12:G15").Value = "O"

# --- 회원가입 (sheet 4): rows 12-21 and 23-24 (row 22 is a merged continuation) ---
$wsSignup = $wb.Worksheets.Item(4)
$wsSignup.Range("E12:E21").Value = "O"
$wsSignup.Range("G12:G21").Value = "O"
$wsSignup.Range("E23:E24").Value = "O"
$wsSignup.Range("G23:G24").Value = "O"

# --- Update selections / active sheet ---
# Set selection on 로그인&로그아웃 and 회원가입 first, then activate
# 메인화면 last so it ends up as the active (tabSelected) sheet, matching
# the target workbook view state.
$wsLogin.Activate()
$wsLogin.Range("G12:G15").Select()

$wsSignup.Activate()
$wsSignup.Range("G12:G29").Select()

$wsMain.Activate()
$wsMain.Range("C32").Select()
